# Append the latest day's profit figure (run on 2025-11-23) as a new row
# at the bottom of the data table, following the existing layout where the
# Date column holds literal text (not a date serial) and Profit holds a
# plain number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCell = $ws.Range("A98")

# Writing the date string directly would cause Excel to auto-convert it
# into a date serial value. Route it through a text formula first so it is
# accepted verbatim, then collapse the formula down to its literal text
# result (Copy + PasteSpecial values) so the cell ends up as plain text,
# matching the style of the other date cells in the column (e.g. A97).
$dateCell.Formula = "=""11/23/2025"""
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)

$ws.Range("B98").Value = 8048.36
